# Hoàn thiện luận cung đại vận
# Adds 5 new rows (6-10) describing the Five Elements (Ngũ hành) relationship
# between the native's Mệnh and the cung đại vận, each paired with its
# corresponding life-outcome description.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New content rows 6-10
$ws.Range("A6").Value = "Ngũ hành bản Mệnh sinh Ngũ hành cung đại vận"
$ws.Range("B6").Value = "Bản thân tổn hao sức khỏe, tiền bạc, công sức. Phải dụng công sức trí lực để hưởng trái ngọt."

$ws.Range("A7").Value = "Ngũ hành bản Mệnh khắc Ngũ hành cung đại vận"
$ws.Range("B7").Value = "Cuộc sống gặp nhiều khó khăn, bản thân phải năng động để có cuộc sống tốt hơn."

$ws.Range("A8").Value = "Ngũ hành cung đại vận sinh Ngũ hành bản Mệnh"
$ws.Range("B8").Value = "Cuộc sống có nhiều thuận lợi, bản thân có nhiều cơ hội và may mắn bất ngờ."

$ws.Range("A9").Value = "Ngũ hành cung đại vận khắc Ngũ hành bản Mệnh"
$ws.Range("B9").Value = "Cuộc bế tắc, bản thân gặp trì trệ, không có nhiều bứt phá."

$ws.Range("A10").Value = "Ngũ hành cung đại vận đồng hành với Ngũ hành bản Mệnh"
$ws.Range("B10").Value = "Cuộc sống thuận lợi, có nhiều cơ hội cho bản thân phát triển."

# Widen column A to fit the longer Vietnamese labels
$ws.Columns.Item(1).ColumnWidth = 59

# Move the active selection back to the top of the sheet
$ws.Range("A1").Select()
$ws.Range("B1").Select()
